$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (row 21) documenting the "Super Queue" problem using Deque.
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Super Queue"
$ws.Range("C21").Value = "Supports finding max element in O(1)"
$ws.Range("H21").Value = "SuperQueue"

# Update selection to the new last cell, mirroring Excel's own post-edit state.
$ws.Range("H21").Select()
